# Append four new ticker rows to the bottom of the "Ticker" column on
# Sheet1, extending the used range from A1:A325 to A1:A329.
# (Daten aktualisiert am 2024-02-23)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A326").Value = "IMX-USD"
$ws.Range("A327").Value = "TAO-USD"
$ws.Range("A328").Value = "MNT-USD"
$ws.Range("A329").Value = "GRT-USD"
